# Update the meanRatingLift (B) and meanQuerySimilarity (C) columns with
# refreshed values, and move the active cell selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$newB = @(-0.1, -0.076, -0.049, -0.022, 0.023, 0.095, 0.244, 0.441, 0.626, 0.711, 0.748)
$newC = @(0.751, 0.75, 0.749, 0.746, 0.74, 0.726, 0.684, 0.61, 0.526, 0.469, 0.415)

for ($i = 0; $i -lt 11; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value = $newB[$i]
    $ws.Cells.Item($row, 3).Value = $newC[$i]
}

$ws.Range("A3").Select()
